$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.222.99"
$ws.Range("E2").Value = "  -1.75%  "
$ws.Range("D3").Value = "2.270.52"
$ws.Range("E3").Value = "  -2.91%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "298.42"
$ws.Range("E5").Value = "  -2.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.16"
$ws.Range("E6").Value = "  -5.84%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -3.23%  "
$ws.Range("E9").Value = "  -3.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.39"
$ws.Range("E10").Value = "  -3.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0791"
$ws.Range("E11").Value = "  -0.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "48.17"
$ws.Range("E12").Value = "  -8.26%  "
$ws.Range("E13").Value = "  -0.33%  "
$ws.Range("E14").Value = "  -2.17%  "
$ws.Range("D15").Value = "2.622.76"
$ws.Range("E15").Value = "  -3.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.55"
$ws.Range("E16").Value = "  -2.36%  "
$ws.Range("D17").Value = "2.265.88"
$ws.Range("E17").Value = "  -2.90%  "
$ws.Range("E18").Value = "  -6.11%  "
$ws.Range("D19").Value = "42.143.87"
$ws.Range("E19").Value = "  -1.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.74"
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("E21").Value = "  -2.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.00"
$ws.Range("E22").Value = "  -2.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.61"
$ws.Range("E23").Value = "  -4.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "234.35"
$ws.Range("E24").Value = "  -0.95%  "
$ws.Range("E25").Value = "  -1.42%  "
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("E27").Value = "  -3.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.91"
$ws.Range("E28").Value = "  -6.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.30"
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.65"
$ws.Range("E30").Value = "  +4.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.25"
$ws.Range("E31").Value = "  -1.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.11"
$ws.Range("E32").Value = "  -1.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.91"
$ws.Range("E34").Value = "  -3.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.55"
$ws.Range("E35").Value = "  -1.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "16.64"
$ws.Range("E36").Value = "  -3.93%  "
$ws.Range("E37").Value = "  -5.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0687"
$ws.Range("E38").Value = "  -4.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.79"
$ws.Range("E39").Value = "  -3.85%  "
$ws.Range("E40").Value = "  -2.93%  "
$ws.Range("E41").Value = "  -2.95%  "
$ws.Range("E42").Value = "  -6.47%  "
$ws.Range("E43").Value = "  -5.03%  "
$ws.Range("D44").Value = "1.961.65"
$ws.Range("E44").Value = "  -2.93%  "
$ws.Range("E45").Value = "  -2.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.56"
$ws.Range("E46").Value = "  -6.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.54"
$ws.Range("E47").Value = "  -6.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.80"
$ws.Range("E48").Value = "  -3.82%  "
$ws.Range("D49").Value = "2.495.21"
$ws.Range("E49").Value = "  -2.45%  "
$ws.Range("E50").Value = "  -6.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.53"
$ws.Range("E51").Value = "  -4.04%  "
